# Prototype features: insert 5 new "Exercises Combo" sub-items under the
# renamed "front-end mechanism (Prototype)" section, and refresh totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 new rows right before the old row 11 ("Jasmine integration"),
#    pushing everything below down by 5 rows. The new rows inherit the
#    indented "sub-item" style from row 10 automatically.
$ws.Rows("11:15").Insert()

# 2) Fill in the 5 newly inserted sub-item rows (11-15).
$ws.Range("C11").Value = "Exercises Combo"
$ws.Range("D11").Value = 4

$ws.Range("C12").Value = "Display Page"
$ws.Range("D12").Value = 4

$ws.Range("C13").Value = "Reduce Page size"
$ws.Range("D13").Value = 4

$ws.Range("C14").Value = "Display review template"
$ws.Range("D14").Value = 4

$ws.Range("C15").Value = "Display Features"
$ws.Range("D15").Value = 4

# Make sure the newly inserted rows truly carry the same indented style as
# the other sub-item rows (C6:C10), in case style didn't fully propagate.
$ws.Range("C11:C15").IndentLevel = 5
$ws.Range("C11:C15").HorizontalAlignment = -4131

# 3) Rename the section header in C5 and replace its SUM formula with the
#    literal rolled-up total (33) now that it has 10 sub-items (D6:D15).
$ws.Range("C5").Value = "front-end mechanism (Prototype)"
$ws.Range("E5").Value = 33

# 4) Old row 11 ("Jasmine integration") is now row 16; it picks up a new
#    D16 contribution (4) alongside its existing E16 total.
$ws.Range("D16").Value = 4

# 5) Refresh the selection to match the edited workbook.
$ws.Range("C6").Select()
